$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Execution Flag column updates
$ws.Range("C2").Value = "YES"
$ws.Range("C4").Value = "NO"
$ws.Range("C6").Value = "YES"

# E3 should carry the Monospace font used elsewhere (e.g. D4)
$ws.Range("D4").Copy()
$ws.Range("E3").PasteSpecial(-4122)

# Clear stray leftover values
$ws.Range("E3").ClearContents()
$ws.Range("F4").ClearContents()

# Move selection to C6
$ws.Range("C6").Select()
